# Update res_line/pl_mw.xlsx values for the 380 kV case (rows 2-25, cols B-O except G/K/M)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.363018046572336
$ws.Cells.Item(2, 3).Value = 0.1319761253530487
$ws.Cells.Item(2, 4).Value = 0.1317583377188782
$ws.Cells.Item(2, 5).Value = 0.1135247653422216
$ws.Cells.Item(2, 6).Value = 1.491495209042576
$ws.Cells.Item(2, 8).Value = 0.07973214163530429
$ws.Cells.Item(2, 9).Value = 0.8675160636340387
$ws.Cells.Item(2, 10).Value = 0.1293826194608663
$ws.Cells.Item(2, 12).Value = 0.3492390935781771
$ws.Cells.Item(2, 14).Value = 1.337447668177361
$ws.Cells.Item(2, 15).Value = 3.774698043901594

# Row 3
$ws.Cells.Item(3, 2).Value = 1.268414032269447
$ws.Cells.Item(3, 3).Value = 0.1183424876848846
$ws.Cells.Item(3, 4).Value = 0.131125836358791
$ws.Cells.Item(3, 5).Value = 0.1140544486402408
$ws.Cells.Item(3, 6).Value = 1.494010841709823
$ws.Cells.Item(3, 8).Value = 0.07973214163530429
$ws.Cells.Item(3, 9).Value = 0.8749475156127282
$ws.Cells.Item(3, 10).Value = 0.1304495154350125
$ws.Cells.Item(3, 12).Value = 0.3418900136483387
$ws.Cells.Item(3, 14).Value = 1.346526100104349
$ws.Cells.Item(3, 15).Value = 3.78417970128092

# Row 4
$ws.Cells.Item(4, 2).Value = 1.210577653479788
$ws.Cells.Item(4, 3).Value = 0.1099183619527651
$ws.Cells.Item(4, 4).Value = 0.1307741754445146
$ws.Cells.Item(4, 5).Value = 0.1144139746642381
$ws.Cells.Item(4, 6).Value = 1.496348933200686
$ws.Cells.Item(4, 8).Value = 0.07973214163530429
$ws.Cells.Item(4, 9).Value = 0.8799707412762743
$ws.Cells.Item(4, 10).Value = 0.1311442677010053
$ws.Cells.Item(4, 12).Value = 0.3375047485304492
$ws.Cells.Item(4, 14).Value = 1.352589935462198
$ws.Cells.Item(4, 15).Value = 3.792141254159873

# Row 5
$ws.Cells.Item(5, 2).Value = 1.18707369094156
$ws.Cells.Item(5, 3).Value = 0.1064723040743445
$ws.Cells.Item(5, 4).Value = 0.1306401417731351
$ws.Cells.Item(5, 5).Value = 0.1145691258371473
$ws.Cells.Item(5, 6).Value = 1.497501328810465
$ws.Cells.Item(5, 8).Value = 0.07973214163530429
$ws.Cells.Item(5, 9).Value = 0.8821334609037415
$ws.Cells.Item(5, 10).Value = 0.1314373718628126
$ws.Cells.Item(5, 12).Value = 0.3357498559434902
$ws.Cells.Item(5, 14).Value = 1.35518427571813
$ws.Cells.Item(5, 15).Value = 3.795923685846674

# Row 6
$ws.Cells.Item(6, 2).Value = 1.183174843045606
$ws.Cells.Item(6, 3).Value = 0.1058992991750358
$ws.Cells.Item(6, 4).Value = 0.1306184467999216
$ws.Cells.Item(6, 5).Value = 0.1145954109848661
$ws.Cells.Item(6, 6).Value = 1.497704741107462
$ws.Cells.Item(6, 8).Value = 0.07973214163530429
$ws.Cells.Item(6, 9).Value = 0.8824995669083187
$ws.Cells.Item(6, 10).Value = 0.1314866450882857
$ws.Cells.Item(6, 12).Value = 0.3354604037075433
$ws.Cells.Item(6, 14).Value = 1.355622514239428
$ws.Cells.Item(6, 15).Value = 3.796584250873821

# Row 7
$ws.Cells.Item(7, 2).Value = 1.210260406003158
$ws.Cells.Item(7, 3).Value = 0.1098719402574488
$ws.Cells.Item(7, 4).Value = 0.1307723302245236
$ws.Cells.Item(7, 5).Value = 0.1144160320760852
$ws.Cells.Item(7, 6).Value = 1.496363666520089
$ws.Cells.Item(7, 8).Value = 0.07973214163530429
$ws.Cells.Item(7, 9).Value = 0.8799994399841466
$ws.Cells.Item(7, 10).Value = 0.1311481801548764
$ws.Cells.Item(7, 12).Value = 0.3374809510700771
$ws.Cells.Item(7, 14).Value = 1.352624424287903
$ws.Cells.Item(7, 15).Value = 3.792190086920783

# Row 8
$ws.Cells.Item(8, 2).Value = 1.330347532172027
$ws.Cells.Item(8, 3).Value = 0.127286350801171
$ws.Cells.Item(8, 4).Value = 0.1315326601370188
$ws.Cells.Item(8, 5).Value = 0.1137002915650669
$ws.Cells.Item(8, 6).Value = 1.49219795441185
$ws.Cells.Item(8, 8).Value = 0.07973214163530429
$ws.Cells.Item(8, 9).Value = 0.8699828794213573
$ws.Cells.Item(8, 10).Value = 0.129742257377738
$ws.Cells.Item(8, 12).Value = 0.3466788552539413
$ws.Cells.Item(8, 14).Value = 1.340476377374102
$ws.Cells.Item(8, 15).Value = 3.777523182387426

# Row 9
$ws.Cells.Item(9, 2).Value = 1.567763831769753
$ws.Cells.Item(9, 3).Value = 0.1610092086419002
$ws.Cells.Item(9, 4).Value = 0.1333131096201114
$ws.Cells.Item(9, 5).Value = 0.1125681463084778
$ws.Cells.Item(9, 6).Value = 1.490323190112733
$ws.Cells.Item(9, 8).Value = 0.07973214163530429
$ws.Cells.Item(9, 9).Value = 0.8539940328816868
$ws.Cells.Item(9, 10).Value = 0.127299538002644
$ws.Cells.Item(9, 12).Value = 0.3657177520953496
$ws.Cells.Item(9, 14).Value = 1.320532702456532
$ws.Cells.Item(9, 15).Value = 3.765745028524663

# Row 10
$ws.Cells.Item(10, 2).Value = 1.743298454633702
$ws.Cells.Item(10, 3).Value = 0.1855193228262806
$ws.Cells.Item(10, 4).Value = 0.1347955671645025
$ws.Cells.Item(10, 5).Value = 0.1119008765462031
$ws.Cells.Item(10, 6).Value = 1.492782960875971
$ws.Cells.Item(10, 8).Value = 0.07973214163530429
$ws.Cells.Item(10, 9).Value = 0.8444766532774324
$ws.Cells.Item(10, 10).Value = 0.125695744995328
$ws.Cells.Item(10, 12).Value = 0.3803091986810045
$ws.Cells.Item(10, 14).Value = 1.308236287807958
$ws.Cells.Item(10, 15).Value = 3.76745868650525

# Row 11
$ws.Cells.Item(11, 2).Value = 1.82338005755571
$ws.Cells.Item(11, 3).Value = 0.1966107164489017
$ws.Cells.Item(11, 4).Value = 0.1355074090039921
$ws.Cells.Item(11, 5).Value = 0.1116328437664276
$ws.Cells.Item(11, 6).Value = 1.494735361813753
$ws.Cells.Item(11, 8).Value = 0.07973214163530429
$ws.Cells.Item(11, 9).Value = 0.8406316217869829
$ws.Cells.Item(11, 10).Value = 0.1250074293809416
$ws.Cells.Item(11, 12).Value = 0.3870768737381951
$ws.Cells.Item(11, 14).Value = 1.303152315004709
$ws.Cells.Item(11, 15).Value = 3.770492679373007

# Row 12
$ws.Cells.Item(12, 2).Value = 1.853736419578695
$ws.Cells.Item(12, 3).Value = 0.2008021969644744
$ws.Cells.Item(12, 4).Value = 0.1357823120924309
$ws.Cells.Item(12, 5).Value = 0.1115364369966017
$ws.Cells.Item(12, 6).Value = 1.495594505083233
$ws.Cells.Item(12, 8).Value = 0.07973214163530429
$ws.Cells.Item(12, 9).Value = 0.8392453196192449
$ws.Cells.Item(12, 10).Value = 0.1247527042709518
$ws.Cells.Item(12, 12).Value = 0.3896581472352949
$ws.Cells.Item(12, 14).Value = 1.301300322304265
$ws.Cells.Item(12, 15).Value = 3.771965926825061

# Row 13
$ws.Cells.Item(13, 2).Value = 1.847197276270322
$ws.Cells.Item(13, 3).Value = 0.1998998717558038
$ws.Cells.Item(13, 4).Value = 0.1357228698052921
$ws.Cells.Item(13, 5).Value = 0.1115569737166471
$ws.Cells.Item(13, 6).Value = 1.495404145143027
$ws.Cells.Item(13, 8).Value = 0.07973214163530429
$ws.Cells.Item(13, 9).Value = 0.8395407824559058
$ws.Cells.Item(13, 10).Value = 0.1248073004845267
$ws.Cells.Item(13, 12).Value = 0.3891014038760119
$ws.Cells.Item(13, 14).Value = 1.301695928146657
$ws.Cells.Item(13, 15).Value = 3.771634209066235

# Row 14
$ws.Cells.Item(14, 2).Value = 1.825876879249279
$ws.Cells.Item(14, 3).Value = 0.1969557253326286
$ws.Cells.Item(14, 4).Value = 0.13552991859153
$ws.Cells.Item(14, 5).Value = 0.1116248103699302
$ws.Cells.Item(14, 6).Value = 1.494803643056997
$ws.Cells.Item(14, 8).Value = 0.07973214163530429
$ws.Cells.Item(14, 9).Value = 0.8405161720730803
$ws.Cells.Item(14, 10).Value = 0.1249863542898844
$ws.Cells.Item(14, 12).Value = 0.3872888669183823
$ws.Cells.Item(14, 14).Value = 1.302998484024968
$ws.Cells.Item(14, 15).Value = 3.770607383188462

# Row 15
$ws.Cells.Item(15, 2).Value = 1.812821527499182
$ws.Cells.Item(15, 3).Value = 0.1951512262584458
$ws.Cells.Item(15, 4).Value = 0.1354124250700863
$ws.Cells.Item(15, 5).Value = 0.1116670248859499
$ws.Cells.Item(15, 6).Value = 1.494451420092858
$ws.Cells.Item(15, 8).Value = 0.07973214163530429
$ws.Cells.Item(15, 9).Value = 0.841122708815611
$ws.Cells.Item(15, 10).Value = 0.1250968013065723
$ws.Cells.Item(15, 12).Value = 0.3861810396870737
$ws.Cells.Item(15, 14).Value = 1.303805866670189
$ws.Cells.Item(15, 15).Value = 3.770020665912085

# Row 16
$ws.Cells.Item(16, 2).Value = 1.738069401397638
$ws.Cells.Item(16, 3).Value = 0.184793282774109
$ws.Cells.Item(16, 4).Value = 0.1347497968018772
$ws.Cells.Item(16, 5).Value = 0.1119191064106921
$ws.Cells.Item(16, 6).Value = 1.492672134729432
$ws.Cells.Item(16, 8).Value = 0.07973214163530429
$ws.Cells.Item(16, 9).Value = 0.8447376891870633
$ws.Cells.Item(16, 10).Value = 0.1257415575885137
$ws.Cells.Item(16, 12).Value = 0.3798695143441506
$ws.Cells.Item(16, 14).Value = 1.308578785692546
$ws.Cells.Item(16, 15).Value = 3.767305786278428

# Row 17
$ws.Cells.Item(17, 2).Value = 1.692268968716462
$ws.Cells.Item(17, 3).Value = 0.1784239288643903
$ws.Cells.Item(17, 4).Value = 0.1343528605154916
$ws.Cells.Item(17, 5).Value = 0.1120828348686214
$ws.Cells.Item(17, 6).Value = 1.49179403894928
$ws.Cells.Item(17, 8).Value = 0.07973214163530429
$ws.Cells.Item(17, 9).Value = 0.8470794911240134
$ws.Cells.Item(17, 10).Value = 0.1261476562662125
$ws.Cells.Item(17, 12).Value = 0.3760307501456595
$ws.Cells.Item(17, 14).Value = 1.311637297256262
$ws.Cells.Item(17, 15).Value = 3.766217822271557

# Row 18
$ws.Cells.Item(18, 2).Value = 1.66594753415427
$ws.Cells.Item(18, 3).Value = 0.1747549615460855
$ws.Cells.Item(18, 4).Value = 0.1341280821780941
$ws.Cells.Item(18, 5).Value = 0.1121803505643673
$ws.Cells.Item(18, 6).Value = 1.491367430873723
$ws.Cells.Item(18, 8).Value = 0.07973214163530429
$ws.Cells.Item(18, 9).Value = 0.8484720308092619
$ws.Cells.Item(18, 10).Value = 0.1263851173288089
$ws.Cells.Item(18, 12).Value = 0.3738350373448469
$ws.Cells.Item(18, 14).Value = 1.313444455934174
$ws.Cells.Item(18, 15).Value = 3.765804243674353

# Row 19
$ws.Cells.Item(19, 2).Value = 1.657039343830775
$ws.Cells.Item(19, 3).Value = 0.1735117764034442
$ws.Cells.Item(19, 4).Value = 0.134052583511064
$ws.Cells.Item(19, 5).Value = 0.1122139423993129
$ws.Cells.Item(19, 6).Value = 1.491236463563368
$ws.Cells.Item(19, 8).Value = 0.07973214163530429
$ws.Cells.Item(19, 9).Value = 0.8489513503086528
$ws.Cells.Item(19, 10).Value = 0.1264661848562634
$ws.Cells.Item(19, 12).Value = 0.3730937150838685
$ws.Cells.Item(19, 14).Value = 1.314064573212008
$ws.Cells.Item(19, 15).Value = 3.765700652108279

# Row 20
$ws.Cells.Item(20, 2).Value = 1.697142264208026
$ws.Cells.Item(20, 3).Value = 0.1791025263851509
$ws.Cells.Item(20, 4).Value = 0.1343947501459297
$ws.Cells.Item(20, 5).Value = 0.1120650597732986
$ws.Cells.Item(20, 6).Value = 1.49187939478486
$ws.Cells.Item(20, 8).Value = 0.07973214163530429
$ws.Cells.Item(20, 9).Value = 0.8468254827309778
$ws.Cells.Item(20, 10).Value = 0.1261040244618385
$ws.Cells.Item(20, 12).Value = 0.3764381274032473
$ws.Cells.Item(20, 14).Value = 1.311306748087517
$ws.Cells.Item(20, 15).Value = 3.766311675950845

# Row 21
$ws.Cells.Item(21, 2).Value = 1.832138364955142
$ws.Cells.Item(21, 3).Value = 0.1978207271721431
$ws.Cells.Item(21, 4).Value = 0.135586448350125
$ws.Cells.Item(21, 5).Value = 0.1116047470406105
$ws.Cells.Item(21, 6).Value = 1.494976773718747
$ws.Cells.Item(21, 8).Value = 0.07973214163530429
$ws.Cells.Item(21, 9).Value = 0.8402277836174576
$ws.Cells.Item(21, 10).Value = 0.1249336011034536
$ws.Cells.Item(21, 12).Value = 0.3878207523001294
$ws.Cells.Item(21, 14).Value = 1.30261390601099
$ws.Cells.Item(21, 15).Value = 3.77090018311651

# Row 22
$ws.Cells.Item(22, 2).Value = 1.920547090840955
$ws.Cells.Item(22, 3).Value = 0.210003996978827
$ws.Cells.Item(22, 4).Value = 0.1363964211121527
$ws.Cells.Item(22, 5).Value = 0.1113335758302743
$ws.Cells.Item(22, 6).Value = 1.497699425404832
$ws.Cells.Item(22, 8).Value = 0.07973214163530429
$ws.Cells.Item(22, 9).Value = 0.8363222613351837
$ws.Cells.Item(22, 10).Value = 0.124203194114159
$ws.Cells.Item(22, 12).Value = 0.3953677127062605
$ws.Cells.Item(22, 14).Value = 1.297359238563651
$ws.Cells.Item(22, 15).Value = 3.77578961132491

# Row 23
$ws.Cells.Item(23, 2).Value = 1.873345751389991
$ws.Cells.Item(23, 3).Value = 0.2035062101198832
$ws.Cells.Item(23, 4).Value = 0.1359612893698312
$ws.Cells.Item(23, 5).Value = 0.1114755951941202
$ws.Cells.Item(23, 6).Value = 1.496182409259319
$ws.Cells.Item(23, 8).Value = 0.07973214163530429
$ws.Cells.Item(23, 9).Value = 0.8383695011217114
$ws.Cells.Item(23, 10).Value = 0.1245898688988714
$ws.Cells.Item(23, 12).Value = 0.391329959269811
$ws.Cells.Item(23, 14).Value = 1.300124750963661
$ws.Cells.Item(23, 15).Value = 3.773006993252608

# Row 24
$ws.Cells.Item(24, 2).Value = 1.694939016506737
$ws.Cells.Item(24, 3).Value = 0.178795754676969
$ws.Cells.Item(24, 4).Value = 0.1343758011707195
$ws.Cells.Item(24, 5).Value = 0.1120730853458198
$ws.Cells.Item(24, 6).Value = 1.491840561754771
$ws.Cells.Item(24, 8).Value = 0.07973214163530429
$ws.Cells.Item(24, 9).Value = 0.8469401759981352
$ws.Cells.Item(24, 10).Value = 0.1261237379731011
$ws.Cells.Item(24, 12).Value = 0.3762539171085564
$ws.Cells.Item(24, 14).Value = 1.311456037411489
$ws.Cells.Item(24, 15).Value = 3.766268584633707

# Row 25
$ws.Cells.Item(25, 2).Value = 1.503337425576547
$ws.Cells.Item(25, 3).Value = 0.1519325806499694
$ws.Cells.Item(25, 4).Value = 0.1328006805201127
$ws.Cells.Item(25, 5).Value = 0.1128454624313697
$ws.Cells.Item(25, 6).Value = 1.490156521343934
$ws.Cells.Item(25, 8).Value = 0.07973214163530429
$ws.Cells.Item(25, 9).Value = 0.8579280601583186
$ws.Cells.Item(25, 10).Value = 0.1279267868464409
$ws.Cells.Item(25, 12).Value = 0.3604607156106283
$ws.Cells.Item(25, 14).Value = 1.325513640109541
$ws.Cells.Item(25, 15).Value = 3.767111731833609

